$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.191.36'
$ws.Range('E2').Value = '  -2.13%  '
$ws.Range('D3').Value = '2.638.22'
$ws.Range('E3').Value = '  -0.27%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '580.44'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.41'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.20%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.650'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.00%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.122'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.41%  '
$ws.Range('E10').Value = '  +0.81%  '
$ws.Range('E11').Value = '  -1.80%  '
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '28.73'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.29%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000187'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.14%  '
$ws.Range('D15').Value = '3.116.97'
$ws.Range('E15').Value = '  -0.03%  '
$ws.Range('D16').Value = '64.053.97'
$ws.Range('E16').Value = '  -2.15%  '
$ws.Range('D17').Value = '2.609.33'
$ws.Range('E17').Value = '  -0.92%  '
$ws.Range('E18').Value = '  -2.67%  '
$ws.Range('E19').Value = '  -0.80%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.56'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '347.49'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.84'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.48%  '
$ws.Range('E24').Value = '  +7.23%  '
$ws.Range('E25').Value = '  -1.87%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.38'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '586.09'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +10.36%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.58'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.45%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.92%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.161'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.70%  '
$ws.Range('E31').Value = '  +0.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.09'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.86%  '
$ws.Range('E33').Value = '  -1.78%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.58'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.28'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.74%  '
$ws.Range('E36').Value = '  -1.48%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.05'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.24%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.998'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.08%  '
$ws.Range('E39').Value = '  +0.17%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '151.70'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('B42').Value = 'OKB'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '42.00'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.50%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '159.35'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.39%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.38'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.96%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.01'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.11%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '23.31'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.61%  '
$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0601'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.04%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.103'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.32%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.635'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.19%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0254'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.35%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.18'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.26%  '
